$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update revised QoQ GVA values for existing rows (column B)
$updates = @{
    7  = -0.407080264463368
    8  = -0.3034247196139308
    23 = 0.03448687557952623
    24 = -0.3400995541928957
    30 = 0.2267802744700904
    31 = -0.04761160248966073
    40 = 1.031257001276817
    41 = -0.3557537420188205
    44 = -0.279615479032766
    45 = 0.01739357871069042
    46 = -0.05522236405459435
    48 = -0.7038072228275515
    49 = -0.2951667760677065
    52 = 0.6867677485477989
    53 = -0.198493249276194
    54 = 0.380581202321879
    55 = -0.6580066892452339
    56 = 0.1112231603001765
    57 = -0.1933116694965538
    59 = 0.2739779345553972
    60 = -0.2626746077356412
    68 = 0.1706552886857859
    69 = -0.1285902175114546
    70 = 0.8762473147740022
    72 = 1.245135280636614
    73 = 0.2066254176843972
    74 = -0.05033923375028097
    75 = -1.280811946119581
    76 = 1.049275222616188
    77 = 0.001054958982941029
    78 = 0.1806871740244445
    80 = 0.05526312921610443
    81 = -0.104378690060197
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}

# Add new row 82: copy formatting from row 81's date cell (A81), then set values
$ws.Cells.Item(81, 1).Copy()
$ws.Cells.Item(82, 1).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(82, 1).Value = 45884
$ws.Cells.Item(82, 2).Value = 0
